# Implement header re-initialization feature test data:
# Insert a new "Cookies!" worksheet before the existing "Annual Report 2022"
# sheet, populate it with a small recipe table, and leave the original
# report sheet selected/active (as the author's workbook view shows).

$wb = $excel.ActiveWorkbook

# Existing (only) sheet -- keep a handle to it before inserting the new one.
$report = $wb.Worksheets.Item(1)
$reportName = $report.Name

# Insert the new sheet ahead of the report sheet so it becomes the first tab.
$cookies = $wb.Worksheets.Add($report)
$cookies.Name = "Cookies!"

$cookies.Range("A1").Value = "Fabulous cookie recipe"
$cookies.Range("A2").Value = "Flour"
$cookies.Range("B2").Value = "2c"
$cookies.Range("A3").Value = "Sugar"
$cookies.Range("B3").Value = "1c"
$cookies.Range("A4").Value = "Butter"
$cookies.Range("B4").Value = ".5lb"
$cookies.Range("A5").Value = "Eggs"
$cookies.Range("B5").Value = 2

# Leave the cursor where the commit's sheet1.xml shows it.
$null = $cookies.Range("D8").Select()

# The Annual Report sheet remains the active tab, with its selection moved.
# (Re-resolving it through $excel.Worksheets, rather than the stale
# reference above, is what actually flips the active tab in this host.)
$excel.Worksheets.Item($reportName).Activate()
$null = $wb.Worksheets.Item($reportName).Range("C5").Select()
